$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, shifting rows 88..202 down to 89..203
$ws.Rows.Item(88).Insert()

# Fill in the new row 88 with the new record's data
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value = "Ñuble"
$ws.Cells.Item(88, 4).Value = 44601
$ws.Cells.Item(88, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = 100112043
$ws.Cells.Item(88, 7).Value = "Pepino ensalada"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 60
$ws.Cells.Item(88, 11).Value = 9500
$ws.Cells.Item(88, 12).Value = 10000
$ws.Cells.Item(88, 13).Value = 9750
$ws.Cells.Item(88, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(88, 15).Value = "Región del Maule"
$ws.Cells.Item(88, 16).Value = 122
$ws.Cells.Item(88, 17).Value = 80
$ws.Cells.Item(88, 18).Value = "Hortaliza"
